$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204175353050232
$ws.Range("B1").Value = 1.688463926315308
$ws.Range("C1").Value = 3.039627075195312
$ws.Range("D1").Value = 3.747160911560059
$ws.Range("E1").Value = 1.43287456035614
